# "Update to slide cover" - add a copyright notice textbox to the bottom of
# the title slide (slide 1), underneath the existing Harvard Extension
# School logo picture.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The real edit produced a shape with id=5 / name="TextBox 4" (the slide
# already has Title=2, Subtitle=3, Picture=1026, so the next auto-assigned
# textbox would be "TextBox 3"). Add + remove a throwaway textbox first so
# the runtime's internal naming counter advances to match ("TextBox 4").
$placeholder = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
$placeholder.Delete()

# Position/size (EMU, expressed in points since AddTextbox takes points):
#   off  x=3348567  y=6421967
#   ext cx=5744633 cy=369332
$textBox = $s.Shapes.AddTextbox(1, 263.6666929133858, 505.6666929133858, 452.3333070866142, 29.081259842519685)
$textBox.Name = "TextBox 4"

$textBox.Fill.Visible = $false

$textBox.TextFrame.WordWrap = -1
$textBox.TextFrame.AutoSize = 1

$textRange = $textBox.TextFrame.TextRange
$textRange.Text = "Copyright 2021, Stephen F Elston. All rights reserved."
$textRange.ParagraphFormat.Alignment = 2

# First 53 characters ("Copyright 2021, Stephen F Elston. All rights reserved")
# are rendered at 11pt; the trailing period keeps the inherited default size.
$mainRun = $textRange.Characters(1, 53)
$mainRun.Font.Size = 11
